$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the swapped/incorrect "Kode Prodi" (I) and "Prodi" (J) sample values -
# these used to reference a vocational "TKJ" class; replace with the
# correct IPA/IPS code + A1/A2 class values.
$ws.Range("I2").Value = "IPA"
$ws.Range("J2").Value = "A1"
$ws.Range("I3").Value = "IPS"
$ws.Range("J3").Value = "A2"

# Add a new "Status" column (M) mirroring the existing bordered header
# style used by the neighbouring K1/L1 cells.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("M1").Value = "Status"

$ws.Range("L1").Copy()
$ws.Range("M2").PasteSpecial(-4122)
$ws.Range("M2").Value = "Aktif"

$ws.Range("L1").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M3").Value = "Non-Aktif"

$excel.CutCopyMode = 0

$ws.Range("J3").Select()
